# Update gh-pages to output generated at 456a3b4
# Applies the same set of cell updates to both the "展览" and "全部类型"
# worksheets (they carry duplicate data in this workbook).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 5 - 江西·ShiningStaR数字互娱嘉年华
    $ws.Range("F5").Value = 8489

    # Row 6 - 上饶·第一届星光次元国风动漫游戏嘉年华暨我和我的cos小伙伴们
    $ws.Range("G6").Value = "不可售"

    # Row 7 - 南昌·第二届龙年动漫展
    $ws.Range("F7").Value = 1498
    $ws.Range("G7").Value = 65

    # Row 8 - 新余·LD02国风动漫嘉年华
    $ws.Range("F8").Value = 174
    $ws.Range("G8").Value = 60

    # Row 11 - 萍乡·AU8春季国漫展
    $ws.Range("F11").Value = 243
    $ws.Range("G11").Value = 55

    # Row 12 - 赣州·漫库书店次元漫展
    $ws.Range("F12").Value = 382
    $ws.Range("G12").Value = 50

    # Row 13 - 南昌·Youth动漫美食嘉年华
    $ws.Range("F13").Value = 239

    # Row 19 - 南昌·ACG CLUB动漫游戏嘉年华
    $ws.Range("F19").Value = 1223

    # Row 20 - 南昌·CM02动漫游戏博览会
    $ws.Range("F20").Value = 168

    # Row 21 - 信丰·端午节UPUP动漫展
    $ws.Range("F21").Value = 76

    # Row 22 - 上饶·ETI动漫节
    $ws.Range("F22").Value = 131

    # Row 23 - 南昌·LY-COSPLAY大会X运动番PRO2.0（非ONLY）
    $ws.Range("F23").Value = 88

    # Row 25 - 九江·第一届异次元动漫嘉年华
    $ws.Range("F25").Value = 66

    # Row 26 - 南昌·第一届异次元动漫嘉年华
    $ws.Range("F26").Value = 104
}
